$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.2025066614151
$ws.Range("B1").Value = 2.143439292907715
$ws.Range("C1").Value = 3.940133571624756
$ws.Range("D1").Value = 3.303717136383057
$ws.Range("E1").Value = 1.118222594261169
